$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$targetA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/387909c2a87116c2f126e42479168bef8726a881/e2e/a.md"
$targetB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/387909c2a87116c2f126e42479168bef8726a881/e2e/b.md"

# ---- Overview sheet: "Ready for handoff" -> "Handed back: in sync with en-US" ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $statusText
$ov.Range("F2").Value = $statusText
$ov.Range("E3").Value = $statusText
$ov.Range("F3").Value = $statusText

# ---- zh-cn sheet: status text, new handback file/date, new "Latest Target File" hyperlink ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

$zh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$zh.Range("K2").Value = "2016-08-16 10:36:45"
$zh.Range("K3").Value = "2016-08-16 10:36:45"

# Recreate hyperlinks so the new "Latest Target File" links (I2/I3) are interleaved
# in the same order Excel would produce them (A2, I2, A3, I3).
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $targetA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
$zh.Hyperlinks.Add($zh.Range("I2"), $targetA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
$zh.Hyperlinks.Add($zh.Range("A3"), $targetB, [System.Type]::Missing, [System.Type]::Missing, "b.md")
$zh.Hyperlinks.Add($zh.Range("I3"), $targetA, [System.Type]::Missing, [System.Type]::Missing, "a.md")

# ---- de-de sheet: status text, new handback file/date, new "Latest Target File" hyperlink ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

$de.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$de.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$de.Range("K2").Value = "2016-08-16 10:36:54"
$de.Range("K3").Value = "2016-08-16 10:36:54"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $targetA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
$de.Hyperlinks.Add($de.Range("I2"), $targetA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
$de.Hyperlinks.Add($de.Range("A3"), $targetB, [System.Type]::Missing, [System.Type]::Missing, "b.md")
$de.Hyperlinks.Add($de.Range("I3"), $targetA, [System.Type]::Missing, [System.Type]::Missing, "a.md")
